# Applies the recorded edit to the deck:
#   1. Slide 6's table switches from the deck's custom "Table_0" style
#      to the built-in table style {2B4DFA71-34E3-41C8-A4AF-A8E069CD92B3}.
#   2. The presentation's design/theme colour scheme switches from the
#      "Integral" palette to the default "Office" palette (the twelve
#      theme colours dk1/lt1/dk2/lt2/accent1-6/hlink/folHlink).

$p = $ppt.ActivePresentation

# -- 1. Table style on slide 6 -------------------------------------------
$slide6 = $p.Slides.Item(6)
for ($i = 1; $i -le $slide6.Shapes.Count; $i++) {
    $shp = $slide6.Shapes.Item($i)
    if ($shp.HasTable) {
        $shp.Table.ApplyStyle("{2B4DFA71-34E3-41C8-A4AF-A8E069CD92B3}")
    }
}

# -- 2. Theme colours: Integral -> Office --------------------------------
function Set-ThemeColorRGB($colorItem, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    # COM RGB() colour values are packed 0x00BBGGRR.
    $colorItem.RGB = ($b * 65536) + ($g * 256) + $r
}

$officeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$themeColors = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 0; $i -lt $officeColors.Length; $i++) {
    Set-ThemeColorRGB $themeColors.Item($i + 1) $officeColors[$i]
}
